$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.315.06'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '2.932.47'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '595.14'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '143.46'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D9').Value = '6.93'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D13').Value = '33.23'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = '3.417.68'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '61.308.66'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '2.934.95'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').Value = '433.68'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('E25').Value = '  -1.74%  '
$ws.Range('D26').Value = '11.72'
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  -3.95%  '
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('D31').Value = '26.65'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +2.20%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '5.61'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  -1.94%  '
$ws.Range('D38').Value = '1.99'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = '8.49'
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('D41').Value = '42.03'
$ws.Range('E41').Value = '  +5.38%  '
$ws.Range('D42').Value = '0.279'
$ws.Range('E42').Value = '  -2.98%  '
$ws.Range('D43').Value = '2.705.12'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '133.32'
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('D46').Value = '363.96'
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D48').Value = '23.58'
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  +0.08%  '
